# Update view counters (column F) on several rows across sheets
# "展览" (Exhibition), "演出" (Performance), "全部类型" (All types)

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F15").Value = 7851
$wsExhibition.Range("F19").Value = 558
$wsExhibition.Range("F34").Value = 6631
$wsExhibition.Range("F37").Value = 217
$wsExhibition.Range("F48").Value = 49

$wsPerformance = $wb.Worksheets.Item("演出")
$wsPerformance.Range("F14").Value = 150

$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Range("F17").Value = 7851
$wsAllTypes.Range("F20").Value = 558
$wsAllTypes.Range("F37").Value = 6631
$wsAllTypes.Range("F39").Value = 217
$wsAllTypes.Range("F47").Value = 49
$wsAllTypes.Range("F49").Value = 150
